$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at row 110. This pushes old rows 110:128 down to 111:129,
# carrying their formulas (with relative refs auto-adjusted) and formatting along.
$ws.Rows.Item(110).Insert()

# Final target data for rows 109:129 (tag id, side length, cross-ratio-0 y/x, cross-ratio-1 y/x).
# Row 110 is the newly-added tag; everything else is the same data as before the insert,
# just shifted down by one row - we re-assert it explicitly to avoid any ambiguity.
$rows = @(
    @{r=109; c=54; d=100; e=30; f=55; g=20; h=50; n=$false; o=$null},
    @{r=110; c=55; d=100; e=30; f=55; g=25; h=55; n=$true;  o=$null},
    @{r=111; c=56; d=100; e=30; f=55; g=20; h=55; n=$false; o=$null},
    @{r=112; c=57; d=100; e=30; f=55; g=25; h=70; n=$true;  o=1},
    @{r=113; c=58; d=100; e=30; f=55; g=25; h=75; n=$false; o=$null},
    @{r=114; c=59; d=100; e=30; f=55; g=20; h=80; n=$false; o=$null},
    @{r=115; c=60; d=100; e=20; f=45; g=20; h=50; n=$false; o=3},
    @{r=116; c=61; d=100; e=20; f=45; g=20; h=55; n=$true;  o=2},
    @{r=117; c=62; d=100; e=20; f=45; g=25; h=70; n=$false; o=$null},
    @{r=118; c=63; d=100; e=20; f=45; g=25; h=75; n=$false; o=$null},
    @{r=119; c=64; d=100; e=20; f=45; g=20; h=80; n=$true;  o=$null},
    @{r=120; c=65; d=100; e=20; f=50; g=20; h=55; n=$false; o=3},
    @{r=121; c=66; d=100; e=20; f=50; g=25; h=70; n=$false; o=$null},
    @{r=122; c=67; d=100; e=20; f=50; g=25; h=75; n=$true;  o=$null},
    @{r=123; c=68; d=100; e=20; f=50; g=20; h=80; n=$false; o=$null},
    @{r=124; c=69; d=100; e=20; f=55; g=25; h=70; n=$true;  o=$null},
    @{r=125; c=70; d=100; e=20; f=55; g=25; h=75; n=$false; o=$null},
    @{r=126; c=71; d=100; e=20; f=55; g=20; h=80; n=$false; o=$null},
    @{r=127; c=72; d=100; e=25; f=70; g=25; h=75; n=$false; o=$null},
    @{r=128; c=73; d=100; e=25; f=70; g=20; h=80; n=$true;  o=$null},
    @{r=129; c=74; d=100; e=25; f=75; g=20; h=80; n=$false; o=$null}
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f
    $ws.Cells.Item($r, 7).Value = $row.g
    $ws.Cells.Item($r, 8).Value = $row.h

    $ws.Cells.Item($r, 9).Formula  = "=E$r*(D$r-F$r)/((D$r-E$r)*F$r)"
    $ws.Cells.Item($r, 10).Formula = "=G$r*(D$r-H$r)/((D$r-G$r)*H$r)"
    $ws.Cells.Item($r, 11).Formula = "=I$r/J$r"
    $ws.Cells.Item($r, 12).Formula = "=ABS(K$r-1)"
    $ws.Cells.Item($r, 13).Formula = "=NOT(OR(K$r<1,L$r<0.05))"

    if ($row.n) {
        $ws.Cells.Item($r, 14).Value = "x"
    } else {
        $ws.Cells.Item($r, 14).ClearContents()
    }

    if ($null -ne $row.o) {
        $ws.Cells.Item($r, 15).Value = $row.o
    } else {
        $ws.Cells.Item($r, 15).ClearContents()
    }
}

# Update the active selection/view to match what Excel shows after editing near the bottom of the table.
$ws.Application.ActiveWindow.ScrollRow = 101
$ws.Range("A128:XFD128").Select()
